$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I, row 4 (header year 2020) - reuse the same format as H4 (bold border row)
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 2020

# Column I, row 5 (Total row) - bold font, right + center alignment, General number format
$ws.Range("I5").Value = 1287.4000000000001
$ws.Range("I5").Font.Bold = $true
$ws.Range("I5").HorizontalAlignment = -4152
$ws.Range("I5").VerticalAlignment = -4108
$ws.Range("I5").NumberFormat = "General"

# Column I, row 6
$ws.Range("I6").Value = 56.6
$ws.Range("I6").NumberFormat = "General"

# Column I, row 7 - stays empty, just gets the same number format/style as row 6
$ws.Range("I7").NumberFormat = "General"

# Column I, row 8
$ws.Range("I8").Value = 2.5
$ws.Range("I8").NumberFormat = "General"

# Column I, row 9
$ws.Range("I9").Value = 9.3000000000000007
$ws.Range("I9").NumberFormat = "General"

# Column I, row 10 - reuse the bottom-border format from H10, then set value + General format
$ws.Range("H10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 0.9
$ws.Range("I10").NumberFormat = "General"

# Match the saved selection/active cell recorded in the workbook
$ws.Range("L9").Select() | Out-Null
